# Updated cryptos list on Wed Jul 26 16:53:23 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns for rows 2-51 with the
# latest scraped coinranking.com figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.322.63'
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").Value = '1.861.57'
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7028'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '238.07'
$ws.Range("D6").ClearFormats()
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07844'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3051'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.75'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +6.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08148'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.34%  '
$ws.Range("D12").Value = '1.878.30'
$ws.Range("E12").Value = '  -0.55%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.220'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7138'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.21'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.27%  '
$ws.Range("D16").Value = '29.390.30'
$ws.Range("E16").Value = '  +0.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.809'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007781'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '238.97'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.18'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.07%  '
$ws.Range("D21").Value = '2.140.91'
$ws.Range("E21").Value = '  +1.69%  '
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.519'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.30'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.899'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.96%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1420'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.07'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.904'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -5.08%  '
$ws.Range("E30").Value = '  -4.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.472'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.298'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.037'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05167'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.99%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.180'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7049'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9995'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.678'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01844'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.695'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.90%  '
$ws.Range("D41").Value = '1.170.00'
$ws.Range("E41").Value = '  +2.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9189'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.39%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.026'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.79%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '71.72'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.33%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4246'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.45%  '
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.80'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.65%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5350'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.752'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.147'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.970'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.04%  '
